$d = $word.ActiveDocument

# Locate the anchor paragraph: the last bullet of the "KEY ACHIEVEMENTS AND
# IMPACT" section ("Expert methodology validated at highest judicial level").
# The two new bullet paragraphs from the diff are inserted right after it,
# and before the following "TECHNICAL SKILLS" Heading2 paragraph.
$anchorRange = $d.Content.Duplicate
$anchorFound = $anchorRange.Find.Execute(
    "Expert methodology validated at highest judicial level",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $anchorFound) {
    throw "Anchor paragraph not found"
}

# Expand to the full paragraph (including its paragraph mark) so we can
# reliably insert a new paragraph right after it.
$fullRange = $d.Range($anchorRange.Start, $anchorRange.End)
$null = $fullRange.Expand(4) # wdParagraph
$null = $fullRange.InsertParagraphAfter()

# Resolve the paragraph index of the anchor so the freshly-created paragraph
# (immediately following it) can be addressed reliably.
$anchorIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Start -eq $fullRange.Start) {
        $anchorIndex = $i
        break
    }
}

# First new bullet paragraph.
$para1 = $d.Paragraphs.Item($anchorIndex + 1)
$para1.Range.Text = "• Breakthrough demographic discovery: Uncovered systematic voter miscoding affecting millions"

# Insert the second new bullet paragraph right after the first one.
$para1.Range.InsertParagraphAfter()
$para2 = $d.Paragraphs.Item($anchorIndex + 2)
$para2.Range.Text = "• 178% accuracy improvement in racial classification algorithms"

# Bold + color the "178%" figure within the second new paragraph to match
# the styling used elsewhere in the document (w:b / w:color 2C3E50).
$figureRange = $para2.Range.Duplicate
$figureFound = $figureRange.Find.Execute(
    "178%", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $figureFound) {
    throw "178% figure not found in new paragraph"
}

$figureRange.Bold = 1
$figureRange.Font.Color = 5258796

Write-Output "Inserted achievement bullets after the judicial-level bullet."
